$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 82, shifting existing rows 82:94 down to 83:95
$ws.Rows.Item(82).Insert()

# Populate the new row 82 with the new weekly price record
$ws.Range("A82").Value = 5
$ws.Range("B82").Value = "Macroferia Regional de Talca"
$ws.Range("C82").Value = "Maule"
$ws.Range("D82").Value = 44476
$ws.Range("E82").Value = 7
$ws.Range("F82").Value = 100112031
$ws.Range("G82").Value = "Poroto verde"
$ws.Range("H82").Value = "Sin especificar"
$ws.Range("I82").Value = "Primera"
$ws.Range("J82").Value = 150
$ws.Range("K82").Value = 42000
$ws.Range("L82").Value = 42000
$ws.Range("M82").Value = 42000
$ws.Range("N82").Value = "$/malla 25 kilos"
$ws.Range("O82").Value = "Región de Arica y Parinacota"
$ws.Range("P82").Value = 1680
$ws.Range("Q82").Value = 25
$ws.Range("R82").Value = "Hortaliza"

Write-Output "done"
